$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.115.56"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.233.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.96%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "297.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "80.88"
$ws.Range("D6").ClearFormats()
$ws.Range("E7").Value = "  -3.67%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.461"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0774"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "28.06"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -9.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.98"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -13.18%  "
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.578.24"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.12"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -7.58%  "
$ws.Range("E16").Value = "  -6.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.241.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.03%  "
$ws.Range("E18").Value = "  -5.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.055.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("E20").Value = "  -5.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.73"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.93"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.22"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.15%  "
$ws.Range("E26").Value = "  -8.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.73"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.93"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "149.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.94%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.79"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -8.12%  "
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0685"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.08%  "
$ws.Range("E37").Value = "  -3.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0968"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.67"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.30%  "
$ws.Range("E41").Value = "  -6.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.65"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.906.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -9.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0255"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.88%  "
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("E48").Value = "  -9.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.441.77"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "87.65"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.15%  "
$ws.Range("E51").Value = "  -7.93%  "
